# CUS15: actualización de servicios, scripts y archivos de cotización
# Adds a "Cotizacion N°" label + number to the quotation header, a stray
# quote-prefixed blank cell at E4, and updates the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COTIZACION")

# --- Header row: new "Cotizacion N°" label (C1), styled like the other
#     section headers (A3 / A9), and its value (D1) ---
$ws.Range("A3").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Cotizacion N°"

$ws.Range("D1").Value = 55555

# --- Stray quote-prefixed empty cell next to "Cliente:" row ---
$ws.Range("E4").Value = "'"

# --- Selection / scroll position update ---
$ws.Range("F8").Select()
